$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.697.04'
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").Value = '2.083.36'
$ws.Range("E3").Value = '  +4.50%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.13'
$ws.Range("E5").Value = '  -2.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  +1.86%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.67'
$ws.Range("E7").Value = '  +6.98%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.387'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.36'
$ws.Range("E10").Value = '  +1.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0765'
$ws.Range("E11").Value = '  +1.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.101'
$ws.Range("E12").Value = '  +3.88%  '

$ws.Range("D13").Value = '2.385.69'
$ws.Range("E13").Value = '  +4.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.58'
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.09'
$ws.Range("E15").Value = '  +1.17%  '

$ws.Range("E16").Value = '  +3.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.27'
$ws.Range("E17").Value = '  +4.38%  '

$ws.Range("D18").Value = '2.073.57'
$ws.Range("E18").Value = '  +3.72%  '

$ws.Range("D19").Value = '37.761.81'
$ws.Range("E19").Value = '  +2.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.08'
$ws.Range("E20").Value = '  +20.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.66'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '224.79'
$ws.Range("E23").Value = '  -1.55%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  +2.94%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  +1.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.11'
$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("E28").Value = '  +2.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.133'
$ws.Range("E29").Value = '  +4.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.41'
$ws.Range("E30").Value = '  +1.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.40'
$ws.Range("E31").Value = '  +6.96%  '

$ws.Range("E32").Value = '  +0.78%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.63'
$ws.Range("E33").Value = '  +13.30%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.50'
$ws.Range("E34").Value = '  +1.47%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0631'
$ws.Range("E35").Value = '  +3.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.46'
$ws.Range("E36").Value = '  +5.52%  '

$ws.Range("E37").Value = '  -0.28%  '

$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.97'
$ws.Range("E38").Value = '  +12.53%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.36'
$ws.Range("E39").Value = '  +0.81%  '

$ws.Range("E40").Value = '  +0.32%  '

$ws.Range("E41").Value = '  -4.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0971'
$ws.Range("E42").Value = '  +9.64%  '

$ws.Range("D43").Value = '1.476.97'
$ws.Range("E43").Value = '  +2.89%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '95.56'
$ws.Range("E44").Value = '  +8.13%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.30'
$ws.Range("E45").Value = '  +22.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0212'
$ws.Range("E46").Value = '  +3.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.52'
$ws.Range("E47").Value = '  +8.42%  '

$ws.Range("E48").Value = '  +0.61%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.37'
$ws.Range("E49").Value = '  +9.42%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.03'
$ws.Range("E50").Value = '  +2.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.94'
$ws.Range("E51").Value = '  +1.83%  '
